$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.523.20"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "2.609.31"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.56"
$ws.Range("E5").Value = "  +2.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.61"
$ws.Range("E6").Value = "  +1.89%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.568"
$ws.Range("E8").Value = "  +0.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.49"
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("E10").Value = "  +1.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.334"
$ws.Range("E11").Value = "  +1.61%  "
$ws.Range("D13").Value = "3.068.68"
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("D14").Value = "59.452.68"
$ws.Range("E14").Value = "  +1.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.66"
$ws.Range("E15").Value = "  +1.36%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000133"
$ws.Range("E16").Value = "  +0.55%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.588.31"
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "341.44"
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("E19").Value = "  +1.80%  "
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.38"
$ws.Range("E21").Value = "  -2.13%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.38"
$ws.Range("E23").Value = "  +2.18%  "
$ws.Range("E24").Value = "  +1.56%  "
$ws.Range("E25").Value = "  -1.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  +3.31%  "
$ws.Range("E28").Value = "  +3.41%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  +5.97%  "
$ws.Range("E31").Value = "  -1.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.87"
$ws.Range("E32").Value = "  +1.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.59"
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.11"
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.837"
$ws.Range("E37").Value = "  +3.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.830"
$ws.Range("E38").Value = "  +0.79%  "
$ws.Range("E39").Value = "  +0.92%  "
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "274.65"
$ws.Range("E41").Value = "  +0.90%  "
$ws.Range("E42").Value = "  +1.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.73"
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0523"
$ws.Range("E45").Value = "  +1.24%  "
$ws.Range("D46").Value = "1.948.74"
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.54"
$ws.Range("E47").Value = "  +3.63%  "
$ws.Range("E48").Value = "  +1.73%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "110.90"
$ws.Range("E50").Value = "  -2.31%  "
$ws.Range("E51").Value = "  +0.49%  "
